$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 11
$ws_ALC.Range("H11").Value = 33.454544
$ws_ALC.Range("I11").Value = 33.454544
$ws_ALC.Range("K11").Value = 33.454544
$ws_ALC.Range("M11").Value = 106.545456

# ALC row 19
$ws_ALC.Range("H19").Value = 1152.9131
$ws_ALC.Range("I19").Value = 1058.1666
$ws_ALC.Range("K19").Value = 1058.1666
$ws_ALC.Range("M19").Value = -883.1666

# ALC row 28
$ws_ALC.Range("H28").Value = 1455.0714
$ws_ALC.Range("I28").Value = 946.5
$ws_ALC.Range("J28").Value = 2726.5
$ws_ALC.Range("K28").Value = 946.5
$ws_ALC.Range("L28").Value = 2726.5
$ws_ALC.Range("M28").Value = -461.5
$ws_ALC.Range("N28").Value = -3696.5

# ALC row 41
$ws_ALC.Range("H41").Value = 974.5
$ws_ALC.Range("I41").Value = 1000
$ws_ALC.Range("J41").Value = 898
$ws_ALC.Range("K41").Value = 1000
$ws_ALC.Range("L41").Value = 898
$ws_ALC.Range("M41").Value = -560
$ws_ALC.Range("N41").Value = -1778

# ALC row 98
$ws_ALC.Range("H98").Value = 1137.6154
$ws_ALC.Range("I98").Value = 565.75
$ws_ALC.Range("K98").Value = 565.75
$ws_ALC.Range("M98").Value = 932.25

# ALC row 122
$ws_ALC.Range("H122").Value = 1137.6154
$ws_ALC.Range("I122").Value = 565.75
$ws_ALC.Range("K122").Value = 1697.25
$ws_ALC.Range("M122").Value = 752.75

# ALC row 132
$ws_ALC.Range("H132").Value = 649.63635
$ws_ALC.Range("I132").Value = 698.3158
$ws_ALC.Range("J132").Value = 341.33334
$ws_ALC.Range("K132").Value = 2094.9474
$ws_ALC.Range("L132").Value = 1024.00002
$ws_ALC.Range("M132").Value = 435.0526
$ws_ALC.Range("N132").Value = -6084.000019999999

# ALC row 137
$ws_ALC.Range("H137").Value = 4507.2
$ws_ALC.Range("I137").Value = 4132.143
$ws_ALC.Range("J137").Value = 6476.25
$ws_ALC.Range("K137").Value = 12396.429
$ws_ALC.Range("L137").Value = 19428.75
$ws_ALC.Range("M137").Value = -9846.429
$ws_ALC.Range("N137").Value = -24528.75

# ALC row 138
$ws_ALC.Range("H138").Value = 1870.8182
$ws_ALC.Range("J138").Value = 3500
$ws_ALC.Range("L138").Value = 10500
$ws_ALC.Range("N138").Value = -20780

# ARM row 61
$ws_ARM.Range("H61").Value = 4676
$ws_ARM.Range("I61").Value = 4399
$ws_ARM.Range("K61").Value = 4399
$ws_ARM.Range("M61").Value = -4187

# ARM row 136
$ws_ARM.Range("H136").Value = 4676
$ws_ARM.Range("I136").Value = 4399
$ws_ARM.Range("K136").Value = 13197
$ws_ARM.Range("M136").Value = -10647

# CRP row 31
$ws_CRP.Range("H31").Value = 2359.6667
$ws_CRP.Range("I31").Value = 2302
$ws_CRP.Range("K31").Value = 2302
$ws_CRP.Range("M31").Value = -2007

# CRP row 34
$ws_CRP.Range("H34").Value = 2359.6667
$ws_CRP.Range("I34").Value = 2302
$ws_CRP.Range("K34").Value = 2302
$ws_CRP.Range("M34").Value = -2100

# CRP row 58
$ws_CRP.Range("H58").Value = 2358.6
$ws_CRP.Range("I58").Value = 2073.375
$ws_CRP.Range("J58").Value = 3499.5
$ws_CRP.Range("K58").Value = 2073.375
$ws_CRP.Range("L58").Value = 3499.5
$ws_CRP.Range("M58").Value = -1870.375
$ws_CRP.Range("N58").Value = -3905.5

# CRP row 62
$ws_CRP.Range("H62").Value = 5278.143
$ws_CRP.Range("I62").Value = 4492
$ws_CRP.Range("J62").Value = 6326.3335
$ws_CRP.Range("K62").Value = 4492
$ws_CRP.Range("L62").Value = 6326.3335
$ws_CRP.Range("M62").Value = -3868
$ws_CRP.Range("N62").Value = -7574.3335

# CRP row 65
$ws_CRP.Range("H65").Value = 5278.143
$ws_CRP.Range("I65").Value = 4492
$ws_CRP.Range("J65").Value = 6326.3335
$ws_CRP.Range("K65").Value = 22460
$ws_CRP.Range("L65").Value = 31631.6675
$ws_CRP.Range("M65").Value = -19340
$ws_CRP.Range("N65").Value = -37871.6675

# CRP row 86
$ws_CRP.Range("H86").Value = 9167.637000000001
$ws_CRP.Range("I86").Value = 8982.777
$ws_CRP.Range("K86").Value = 8982.777
$ws_CRP.Range("M86").Value = -7859.777

# CRP row 89
$ws_CRP.Range("H89").Value = 9167.637000000001
$ws_CRP.Range("I89").Value = 8982.777
$ws_CRP.Range("K89").Value = 44913.885
$ws_CRP.Range("M89").Value = -39297.885

# CRP row 94
$ws_CRP.Range("H94").Value = 1328.3334
$ws_CRP.Range("I94").Value = 1328.3334
$ws_CRP.Range("J94").Value = 0
$ws_CRP.Range("K94").Value = 1328.3334
$ws_CRP.Range("L94").Value = 0
$ws_CRP.Range("M94").Value = -877.3334
$ws_CRP.Range("N94").ClearContents()

# CRP row 129
$ws_CRP.Range("H129").Value = 0
$ws_CRP.Range("J129").Value = 0
$ws_CRP.Range("L129").Value = 0
$ws_CRP.Range("N129").ClearContents()

# CRP row 132
$ws_CRP.Range("H132").Value = 1287.7
$ws_CRP.Range("I132").Value = 1211.1428
$ws_CRP.Range("K132").Value = 3633.4284
$ws_CRP.Range("M132").Value = -1103.4284

# CRP row 134
$ws_CRP.Range("H134").Value = 1499.5
$ws_CRP.Range("I134").Value = 1473
$ws_CRP.Range("K134").Value = 4419
$ws_CRP.Range("M134").Value = -1884

# CRP row 136
$ws_CRP.Range("H136").Value = 2358.6
$ws_CRP.Range("I136").Value = 2073.375
$ws_CRP.Range("J136").Value = 3499.5
$ws_CRP.Range("K136").Value = 6220.125
$ws_CRP.Range("L136").Value = 10498.5
$ws_CRP.Range("M136").Value = -3670.125
$ws_CRP.Range("N136").Value = -15598.5

# CUL row 62
$ws_CUL.Range("H62").Value = 10533.5
$ws_CUL.Range("J62").Value = 10511.333
$ws_CUL.Range("L62").Value = 31533.999
$ws_CUL.Range("N62").Value = -32905.999

# CUL row 65
$ws_CUL.Range("H65").Value = 10533.5
$ws_CUL.Range("J65").Value = 10511.333
$ws_CUL.Range("L65").Value = 94601.997
$ws_CUL.Range("N65").Value = -101465.997

# CUL row 98
$ws_CUL.Range("H98").Value = 1242.5
$ws_CUL.Range("J98").Value = 990
$ws_CUL.Range("L98").Value = 2970
$ws_CUL.Range("N98").Value = -5966

# CUL row 103
$ws_CUL.Range("H103").Value = 999.2
$ws_CUL.Range("I103").Value = 24.5
$ws_CUL.Range("J103").Value = 1649
$ws_CUL.Range("K103").Value = 73.5
$ws_CUL.Range("L103").Value = 4947
$ws_CUL.Range("M103").Value = 805.5
$ws_CUL.Range("N103").Value = -6705

# CUL row 131
$ws_CUL.Range("H131").Value = 879
$ws_CUL.Range("I131").Value = 879
$ws_CUL.Range("K131").Value = 2637
$ws_CUL.Range("M131").Value = 2403

# GSM row 9
$ws_GSM.Range("H9").Value = 1981.1
$ws_GSM.Range("I9").Value = 780.8
$ws_GSM.Range("J9").Value = 3181.4
$ws_GSM.Range("K9").Value = 780.8
$ws_GSM.Range("L9").Value = 3181.4
$ws_GSM.Range("M9").Value = -610.8
$ws_GSM.Range("N9").Value = -3521.4

# LTW row 54
$ws_LTW.Range("H54").Value = 10000
$ws_LTW.Range("J54").Value = 10000
$ws_LTW.Range("L54").Value = 10000
$ws_LTW.Range("N54").Value = -11288

# WVR row 107
$ws_WVR.Range("H107").Value = 1089.6666
$ws_WVR.Range("I107").Value = 1262.8334
$ws_WVR.Range("J107").Value = 743.3333
$ws_WVR.Range("K107").Value = 3788.5002
$ws_WVR.Range("L107").Value = 2229.9999
$ws_WVR.Range("M107").Value = -1868.5002
$ws_WVR.Range("N107").Value = -6069.9999
